$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '51.835.85'
$ws.Range('E2').Value = '  -0.14%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.829.00'
$ws.Range('E3').Value = '  +1.15%  '

$ws.Range('E4').Value = '  +0.03%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '350.05'
$ws.Range('E5').Value = '  -1.21%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '113.26'
$ws.Range('E6').Value = '  +3.67%  '

$ws.Range('E7').Value = '  +0.46%  '

$ws.Range('E8').Value = '  -0.04%  '

$ws.Range('E9').Value = '  +3.02%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '40.23'
$ws.Range('E10').Value = '  +0.57%  '

$ws.Range('E11').Value = '  -1.01%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0849'
$ws.Range('E12').Value = '  +1.04%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '20.05'
$ws.Range('E13').Value = '  -0.73%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.78'
$ws.Range('E14').Value = '  +1.40%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.271.41'
$ws.Range('E15').Value = '  +1.16%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.985'
$ws.Range('E16').Value = '  +5.82%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.816.45'
$ws.Range('E17').Value = '  +0.30%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '51.804.48'
$ws.Range('E18').Value = '  +0.01%  '

$ws.Range('E19').Value = '  +9.27%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.64'
$ws.Range('E20').Value = '  -1.03%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.42'
$ws.Range('E21').Value = '  +1.56%  '

$ws.Range('E22').Value = '  +0.46%  '

$ws.Range('E23').Value = '  +0.33%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '268.87'
$ws.Range('E24').Value = '  +0.39%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.75'
$ws.Range('E25').Value = '  +0.71%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '26.26'
$ws.Range('E26').Value = '  +0.09%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.00'
$ws.Range('E27').Value = '  -0.13%  '

$ws.Range('E28').Value = '  +1.23%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '39.29'
$ws.Range('E29').Value = '  +6.83%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '10.55'
$ws.Range('E30').Value = '  +2.58%  '

$ws.Range('E31').Value = '  +16.65%  '

$ws.Range('E32').Value = '  +1.52%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '52.74'
$ws.Range('E33').Value = '  +1.30%  '

$ws.Range('B34').Value = 'RenderToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.66'
$ws.Range('E34').Value = '  +1.59%  '

$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0894'
$ws.Range('E35').Value = '  +7.60%  '

$ws.Range('E36').Value = '  -1.21%  '

$ws.Range('E37').Value = '  +0.04%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '18.85'
$ws.Range('E38').Value = '  +1.17%  '

$ws.Range('E39').Value = '  +1.74%  '

$ws.Range('E40').Value = '  +1.89%  '

$ws.Range('E41').Value = '  +0.96%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.52'
$ws.Range('E42').Value = '  -1.85%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '122.64'
$ws.Range('E43').Value = '  +1.35%  '

$ws.Range('E44').Value = '  +1.48%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '22.15'
$ws.Range('E45').Value = '  -0.44%  '

$ws.Range('B46').Value = 'NEARProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.52'
$ws.Range('E46').Value = '  +6.60%  '

$ws.Range('B47').Value = 'ApeXProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.51'
$ws.Range('E47').Value = '  +8.63%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.174.54'
$ws.Range('E48').Value = '  +1.84%  '

$ws.Range('E49').Value = '  +22.28%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.948'
$ws.Range('E50').Value = '  +3.78%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '5.44'
$ws.Range('E51').Value = '  -0.52%  '
